$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 658.8
$ws.Range("I9").Value = 338.33334
$ws.Range("J9").Value = 1139.5
$ws.Range("K9").Value = 338.33334
$ws.Range("L9").Value = 1139.5
$ws.Range("M9").Value = -169.33334
$ws.Range("N9").Value = -1477.5
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("M12").Value = -130
$ws.Range("H21").Value = 28000
$ws.Range("I21").Value = 28000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -27532
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 28000
$ws.Range("I23").Value = 28000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 28000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -27766
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 1659.3334
$ws.Range("J29").Value = 2389
$ws.Range("L29").Value = 7167
$ws.Range("N29").Value = -7729
$ws.Range("H32").Value = 9107.691999999999
$ws.Range("J32").Value = 10740.909
$ws.Range("L32").Value = 10740.909
$ws.Range("N32").Value = -11392.909
$ws.Range("H38").Value = 412.2
$ws.Range("I38").Value = 155.92857
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 467.78571
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = -95.78570999999999
$ws.Range("N38").Value = -12744
$ws.Range("H58").Value = 2529.8
$ws.Range("I58").Value = 570
$ws.Range("J58").Value = 5469.5
$ws.Range("K58").Value = 1710
$ws.Range("L58").Value = 16408.5
$ws.Range("M58").Value = -1560
$ws.Range("N58").Value = -16708.5
$ws.Range("H98").Value = 10255.692
$ws.Range("I98").Value = 10943.667
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 10943.667
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -9445.666999999999
$ws.Range("N98").Value = -4996
$ws.Range("H99").Value = 1235.2941
$ws.Range("J99").Value = 2498
$ws.Range("L99").Value = 7494
$ws.Range("N99").Value = -10490
$ws.Range("H106").Value = 5442.8438
$ws.Range("I106").Value = 6280.5557
$ws.Range("J106").Value = 919.2
$ws.Range("K106").Value = 6280.5557
$ws.Range("L106").Value = 919.2
$ws.Range("M106").Value = -5649.5557
$ws.Range("N106").Value = -2181.2
$ws.Range("H112").Value = 1683.3462
$ws.Range("J112").Value = 1683.3462
$ws.Range("L112").Value = 5050.0386
$ws.Range("N112").Value = -7266.0386
$ws.Range("H122").Value = 10255.692
$ws.Range("I122").Value = 10943.667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 32831.001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -30381.001
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 7760410.5
$ws.Range("I132").Value = 9529832
$ws.Range("J132").Value = 19194.375
$ws.Range("K132").Value = 28589496
$ws.Range("L132").Value = 57583.125
$ws.Range("M132").Value = -28586966
$ws.Range("N132").Value = -62643.125
$ws.Range("H137").Value = 1511.4333
$ws.Range("I137").Value = 1396
$ws.Range("J137").Value = 1612.4375
$ws.Range("K137").Value = 4188
$ws.Range("L137").Value = 4837.3125
$ws.Range("M137").Value = -1638
$ws.Range("N137").Value = -9937.3125
$ws.Range("H138").Value = 3712.1428
$ws.Range("I138").Value = 3715.75
$ws.Range("J138").Value = 3711.9624
$ws.Range("K138").Value = 11147.25
$ws.Range("L138").Value = 11135.8872
$ws.Range("M138").Value = -6007.25
$ws.Range("N138").Value = -21415.8872

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2192.8125
$ws.Range("I2").Value = 1464.75
$ws.Range("K2").Value = 1464.75
$ws.Range("M2").Value = -1351.75
$ws.Range("H32").Value = 25638.547
$ws.Range("I32").Value = 16227.467
$ws.Range("J32").Value = 49166.25
$ws.Range("K32").Value = 16227.467
$ws.Range("L32").Value = 49166.25
$ws.Range("M32").Value = -15940.467
$ws.Range("N32").Value = -49740.25
$ws.Range("H61").Value = 35715560
$ws.Range("I61").Value = 41667828
$ws.Range("J61").Value = 1951
$ws.Range("K61").Value = 41667828
$ws.Range("L61").Value = 1951
$ws.Range("M61").Value = -41667616
$ws.Range("N61").Value = -2375
$ws.Range("H74").Value = 1911
$ws.Range("I74").Value = 930.1667
$ws.Range("J74").Value = 2530.4736
$ws.Range("K74").Value = 930.1667
$ws.Range("L74").Value = 2530.4736
$ws.Range("M74").Value = -56.16669999999999
$ws.Range("N74").Value = -4278.473599999999
$ws.Range("H77").Value = 1911
$ws.Range("I77").Value = 930.1667
$ws.Range("J77").Value = 2530.4736
$ws.Range("K77").Value = 4650.8335
$ws.Range("L77").Value = 12652.368
$ws.Range("M77").Value = -282.8334999999997
$ws.Range("N77").Value = -21388.368
$ws.Range("H97").Value = 6856.5
$ws.Range("I97").Value = 657.4286
$ws.Range("K97").Value = 657.4286
$ws.Range("M97").Value = -161.4286
$ws.Range("H116").Value = 2192.8125
$ws.Range("I116").Value = 1464.75
$ws.Range("K116").Value = 1464.75
$ws.Range("M116").Value = 829.25
$ws.Range("H122").Value = 4193
$ws.Range("I122").Value = 3864.6
$ws.Range("J122").Value = 5014
$ws.Range("K122").Value = 11593.8
$ws.Range("L122").Value = 15042
$ws.Range("M122").Value = -9143.799999999999
$ws.Range("N122").Value = -19942
$ws.Range("H132").Value = 2184.5715
$ws.Range("I132").Value = 1445.5
$ws.Range("J132").Value = 4549.6
$ws.Range("K132").Value = 4336.5
$ws.Range("L132").Value = 13648.8
$ws.Range("M132").Value = -1806.5
$ws.Range("N132").Value = -18708.8
$ws.Range("H134").Value = 37400
$ws.Range("J134").Value = 37400
$ws.Range("L134").Value = 37400
$ws.Range("N134").Value = -47540
$ws.Range("H136").Value = 35715560
$ws.Range("I136").Value = 41667828
$ws.Range("J136").Value = 1951
$ws.Range("K136").Value = 125003484
$ws.Range("L136").Value = 5853
$ws.Range("M136").Value = -125000934
$ws.Range("N136").Value = -10953

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2192.8125
$ws.Range("I3").Value = 1464.75
$ws.Range("K3").Value = 1464.75
$ws.Range("M3").Value = -1350.75
$ws.Range("H22").Value = 600
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -946
$ws.Range("H36").Value = 1297.8334
$ws.Range("I36").Value = 557.4
$ws.Range("K36").Value = 557.4
$ws.Range("M36").Value = -23.39999999999998
$ws.Range("H94").Value = 13158393
$ws.Range("I94").Value = 13158393
$ws.Range("K94").Value = 13158393
$ws.Range("M94").Value = -13157942
$ws.Range("H105").Value = 34484320
$ws.Range("I105").Value = 43479668
$ws.Range("J105").Value = 2162.8333
$ws.Range("K105").Value = 43479668
$ws.Range("L105").Value = 2162.8333
$ws.Range("M105").Value = -43477921
$ws.Range("N105").Value = -5656.8333
$ws.Range("H107").Value = 1400
$ws.Range("I107").Value = 1400
$ws.Range("K107").Value = 1400
$ws.Range("M107").Value = 520
$ws.Range("H134").Value = 3388.4773
$ws.Range("I134").Value = 831.4583
$ws.Range("K134").Value = 2494.3749
$ws.Range("M134").Value = 40.6251000000002

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2201.6
$ws.Range("J4").Value = 2201.6
$ws.Range("L4").Value = 2201.6
$ws.Range("N4").Value = -2425.6
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 200
$ws.Range("M17").Value = -26
$ws.Range("H19").Value = 107.5
$ws.Range("I19").Value = 107.5
$ws.Range("K19").Value = 107.5
$ws.Range("M19").Value = 62.5
$ws.Range("H20").Value = 48864.5
$ws.Range("J20").Value = 48864.5
$ws.Range("L20").Value = 48864.5
$ws.Range("N20").Value = -49336.5
$ws.Range("H24").Value = 107.5
$ws.Range("I24").Value = 107.5
$ws.Range("K24").Value = 107.5
$ws.Range("M24").Value = 62.5
$ws.Range("H30").Value = 48864.5
$ws.Range("J30").Value = 48864.5
$ws.Range("L30").Value = 48864.5
$ws.Range("N30").Value = -49046.5
$ws.Range("H31").Value = 1239.5867
$ws.Range("I31").Value = 1260.4203
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 1260.4203
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -965.4203
$ws.Range("N31").Value = -1590
$ws.Range("H34").Value = 1239.5867
$ws.Range("I34").Value = 1260.4203
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1260.4203
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = -1058.4203
$ws.Range("N34").Value = -1404
$ws.Range("H62").Value = 11113578
$ws.Range("I62").Value = 2449.9167
$ws.Range("K62").Value = 2449.9167
$ws.Range("M62").Value = -1825.9167
$ws.Range("H65").Value = 11113578
$ws.Range("I65").Value = 2449.9167
$ws.Range("K65").Value = 12249.5835
$ws.Range("M65").Value = -9129.583500000001
$ws.Range("H99").Value = 5264757
$ws.Range("I99").Value = 13158893
$ws.Range("K99").Value = 13158893
$ws.Range("M99").Value = -13157395
$ws.Range("H126").Value = 5264757
$ws.Range("I126").Value = 13158893
$ws.Range("K126").Value = 39476679
$ws.Range("M126").Value = -39474209
$ws.Range("H128").Value = 48864.5
$ws.Range("J128").Value = 48864.5
$ws.Range("L128").Value = 48864.5
$ws.Range("N128").Value = -58824.5
$ws.Range("H131").Value = 10577.2
$ws.Range("H132").Value = 2344.6
$ws.Range("I132").Value = 2043.3334
$ws.Range("J132").Value = 3549.6667
$ws.Range("K132").Value = 6130.0002
$ws.Range("L132").Value = 10649.0001
$ws.Range("M132").Value = -3600.0002
$ws.Range("N132").Value = -15709.0001
$ws.Range("H134").Value = 10001484
$ws.Range("I134").Value = 1370.5676
$ws.Range("J134").Value = 38463344
$ws.Range("K134").Value = 4111.7028
$ws.Range("L134").Value = 115390032
$ws.Range("M134").Value = -1576.7028
$ws.Range("N134").Value = -115395102

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H113").Value = 684.9666999999999
$ws.Range("I113").Value = 583.75
$ws.Range("K113").Value = 1751.25
$ws.Range("M113").Value = 418.75
$ws.Range("H122").Value = 1110.0385
$ws.Range("J122").Value = 1521.9286
$ws.Range("L122").Value = 13697.3574
$ws.Range("N122").Value = -18597.3574
$ws.Range("H127").Value = 2448.3333
$ws.Range("J127").Value = 2448.3333
$ws.Range("L127").Value = 7344.999899999999
$ws.Range("N127").Value = -17264.9999
$ws.Range("H131").Value = 15897790
$ws.Range("I131").Value = 125000490
$ws.Range("J131").Value = 28306.328
$ws.Range("K131").Value = 375001470
$ws.Range("L131").Value = 84918.984
$ws.Range("M131").Value = -374996430
$ws.Range("N131").Value = -94998.984

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3171.4092
$ws.Range("I102").Value = 2477.4167
$ws.Range("J102").Value = 4004.2
$ws.Range("K102").Value = 2477.4167
$ws.Range("L102").Value = 4004.2
$ws.Range("M102").Value = -855.4167000000002
$ws.Range("N102").Value = -7248.2
$ws.Range("H122").Value = 1779
$ws.Range("I122").Value = 1682.8
$ws.Range("J122").Value = 2099.6667
$ws.Range("K122").Value = 5048.4
$ws.Range("L122").Value = 6299.000100000001
$ws.Range("M122").Value = -2598.4
$ws.Range("N122").Value = -11199.0001
$ws.Range("H132").Value = 10017.3125
$ws.Range("I132").Value = 13485.1
$ws.Range("K132").Value = 40455.3
$ws.Range("M132").Value = -37925.3

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3998.75
$ws.Range("I40").Value = 3998.3333
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3998.3333
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -3862.3333
$ws.Range("N40").Value = -4272
$ws.Range("H93").Value = 617.75
$ws.Range("I93").Value = 486.17648
$ws.Range("K93").Value = 486.17648
$ws.Range("M93").Value = 761.8235199999999
$ws.Range("H95").Value = 14750
$ws.Range("J95").Value = 14750
$ws.Range("L95").Value = 14750
$ws.Range("N95").Value = -20242
$ws.Range("H100").Value = 1461.7693
$ws.Range("I100").Value = 1363.909
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1363.909
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -822.9090000000001
$ws.Range("N100").Value = -3082
$ws.Range("H122").Value = 40478244
$ws.Range("I122").Value = 56668336
$ws.Range("K122").Value = 170005008
$ws.Range("M122").Value = -170002558
$ws.Range("H132").Value = 4999.4
$ws.Range("I132").Value = 4333
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 12999
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -10469
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 35860
$ws.Range("J134").Value = 35860
$ws.Range("L134").Value = 35860
$ws.Range("N134").Value = -46000
$ws.Range("H136").Value = 2272
$ws.Range("I136").Value = 1468
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 4404
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -1854
$ws.Range("N136").Value = -13725

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6800
$ws.Range("J45").Value = 6800
$ws.Range("L45").Value = 6800
$ws.Range("N45").Value = -7782
$ws.Range("H97").Value = 19142.25
$ws.Range("J97").Value = 19142.25
$ws.Range("L97").Value = 19142.25
$ws.Range("N97").Value = -21124.25
$ws.Range("H100").Value = 358.15384
$ws.Range("I100").Value = 362.75
$ws.Range("J100").Value = 303
$ws.Range("K100").Value = 725.5
$ws.Range("L100").Value = 606
$ws.Range("M100").Value = -184.5
$ws.Range("N100").Value = -1688
$ws.Range("H107").Value = 430.45456
$ws.Range("I107").Value = 376.1111
$ws.Range("K107").Value = 1128.3333
$ws.Range("M107").Value = 791.6667
$ws.Range("H122").Value = 10871560
$ws.Range("I122").Value = 10871560
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 32614680
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -32612230
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 1787.5
$ws.Range("I136").Value = 1060
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3180
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -630
$ws.Range("N136").Value = -14100
